$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")

# Apply check-apl number corrections (values replaced row by row); the apl (applicator)
# counts are now tracked in a separate equipment file, so most per-apparatus counters
# collapse to 0 here.
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = 9
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 23
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("Y3").Value = 0
$ws.Range("AA3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = 6
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 6
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 23
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 0
$ws.Range("V18").Value = 0
$ws.Range("W18").Value = 0
$ws.Range("X18").Value = 0
$ws.Range("Y18").Value = 0
$ws.Range("Z18").Value = 0
$ws.Range("AA18").Value = 0

# Give the apl total cell (Q17) a dedicated number style (thousands separator, right aligned)
$ws.Range("Q17").NumberFormat = "#,##0"

# Update selection to reflect the last active cell used while editing
$ws.Range("P23").Select()
